# Update docs: add a new "-sumer / consumer" word-family row to the
# "root-v" sheet (inserted just before the "-tendre" group, i.e. at
# worksheet row 27) and leave the workbook with that sheet active /
# selected (it was previously on "prefix-shift").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("root-v")

# Insert a new blank row above the current row 27 ("-tendre" group),
# pushing every following row down by one (rows 27-37 -> 28-38).
$ws.Range("A27:E27").EntireRow.Insert()

# Populate the freshly inserted row with the new word family.
$ws.Range("A27").Value = "-sumer"
$ws.Range("B27").Value = "-sume"
$ws.Range("C27").Value = "consumer"
$ws.Range("D27").Value = "consume"
$ws.Range("E27").Value = "+er"

# Make "root-v" the active sheet/tab, with A28 selected (matches the
# sheet's new sheetView selection in the saved workbook).
$ws.Activate()
$ws.Range("A28").Select()
